$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: ODE now uses VM1, VDM1, VM3, VDM3 (instead of VD1/VD3).
# Add the four new parameter rows (53-56) with their labels/values,
# mirroring the formatting used for the other "scaling factor" rows
# just above them (51-52): plain "Normal" style cells with black font.

$ws.Range("E53").Value = "VM1"
$ws.Range("F53").Value = 3.2837957690000001

$ws.Range("E54").Value = "VM3"
$ws.Range("F54").Value = 0.5

$ws.Range("E55").Value = "VDM1"
$ws.Range("F55").Value = 3.2837957690000001

$ws.Range("E56").Value = "VDM3"
$ws.Range("F56").Value = 0.5

$newRange = $ws.Range("E53:F56")
$newRange.Style = "Normal"
$newRange.Font.Color = 0

# Reflect the new selection / scroll position left behind by the edit.
$ws.Range("A19").Select() | Out-Null
$newRange.Select() | Out-Null
